$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '43.081.83'
$ws.Range("E2").Value = '  +4.89%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.278.77'
$ws.Range("E3").Value = '  +4.85%  '

$ws.Range("E4").Value = '  +0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '250.39'
$ws.Range("E5").Value = '  +1.35%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.638'
$ws.Range("E6").Value = '  +3.65%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '72.29'
$ws.Range("E7").Value = '  +8.80%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.671'
$ws.Range("E8").Value = '  +18.46%  '

$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.00'
$ws.Range("E9").Value = '  -0.09%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.39'
$ws.Range("E10").Value = '  +10.94%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '59.76'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0967'
$ws.Range("E12").Value = '  +4.16%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '7.49'
$ws.Range("E13").Value = '  +9.09%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.104'
$ws.Range("E14").Value = '  +0.84%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.616.92'
$ws.Range("E15").Value = '  +4.76%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '15.03'
$ws.Range("E16").Value = '  +5.47%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.889'
$ws.Range("E17").Value = '  +4.00%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.262.37'
$ws.Range("E18").Value = '  +4.96%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '43.041.30'
$ws.Range("E19").Value = '  +5.09%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0000101'
$ws.Range("E20").Value = '  +7.22%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.35'
$ws.Range("E21").Value = '  +4.35%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '73.39'
$ws.Range("E22").Value = '  +2.76%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '236.88'
$ws.Range("E23").Value = '  +3.08%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.11'
$ws.Range("E24").Value = '  +2.10%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '3.95'
$ws.Range("E25").Value = '  +7.37%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '11.55'
$ws.Range("E26").Value = '  +1.77%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.00'
$ws.Range("E27").Value = '  +0.12%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.44'
$ws.Range("E28").Value = '  +0.48%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '3.67'
$ws.Range("E29").Value = '  -1.27%  '

$ws.Range("E30").Value = '  +5.45%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '167.99'
$ws.Range("E31").Value = '  -0.30%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.12'
$ws.Range("E32").Value = '  +4.55%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.60'
$ws.Range("E33").Value = '  +17.05%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.128'
$ws.Range("E34").Value = '  +5.08%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0814'
$ws.Range("E35").Value = '  +8.52%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '31.58'
$ws.Range("E36").Value = '  +29.77%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.126'
$ws.Range("E37").Value = '  +4.31%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.48'
$ws.Range("E38").Value = '  +11.64%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '4.79'
$ws.Range("E39").Value = '  +4.82%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0320'
$ws.Range("E40").Value = '  +5.70%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.34'
$ws.Range("E41").Value = '  +6.85%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '12.69'
$ws.Range("E42").Value = '  +12.45%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.87'
$ws.Range("E43").Value = '  +7.53%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '9.27'
$ws.Range("E44").Value = '  +9.61%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '62.79'
$ws.Range("E45").Value = '  +3.86%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.206'
$ws.Range("E46").Value = '  +7.13%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.87'
$ws.Range("E47").Value = '  +0.22%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.103'
$ws.Range("E48").Value = '  +3.77%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  +3.02%  '

$ws.Range("E51").Value = '  +4.39%  '
